$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, pushing existing rows 11-69 down to 12-70.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly data point.
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 44819
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = 100112035
$ws.Cells.Item(11, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 52
$ws.Cells.Item(11, 11).Value = 14000
$ws.Cells.Item(11, 12).Value = 15000
$ws.Cells.Item(11, 13).Value = 14500
$ws.Cells.Item(11, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(11, 15).Value = "Hijuelas"
$ws.Cells.Item(11, 16).Value = 967
$ws.Cells.Item(11, 17).Value = 15
$ws.Cells.Item(11, 18).Value = "Hortaliza"
